$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2654603333333334
$ws.Range("H2").Value = 0.796381
$ws.Range("I2").Value = 0.04079010536687974
$ws.Range("J2").Value = 0.04079010536687975
$ws.Range("M2").Value = 1.485259333333333
$ws.Range("N2").Value = 4.455778
$ws.Range("O2").Value = 0.3057455162066235
$ws.Range("P2").Value = 0.3057455162066235
$ws.Range("Q2").Value = 0.3942774377131112
$ws.Range("R2").Value = 3.548496939418
$ws.Range("S2").Value = 0.01247139182151921
$ws.Range("T2").Value = 0.01247139182151921
$ws.Range("G3").Value = 0.2654603333333334
$ws.Range("H3").Value = 0.796381
$ws.Range("I3").Value = 0.04079010536687974
$ws.Range("J3").Value = 0.04079010536687975
$ws.Range("O3").Value = 0.2805555239151429
$ws.Range("P3").Value = 0.2805555239151429
$ws.Range("Q3").Value = 0.3617934106702223
$ws.Range("R3").Value = 3.256140696032001
$ws.Range("S3").Value = 0.01144388938175883
$ws.Range("T3").Value = 0.01144388938175883
$ws.Range("G4").Value = 0.2654603333333334
$ws.Range("H4").Value = 0.796381
$ws.Range("I4").Value = 0.04079010536687974
$ws.Range("J4").Value = 0.04079010536687975
$ws.Range("O4").Value = 0.4136989598782336
$ws.Range("P4").Value = 0.4136989598782336
$ws.Range("Q4").Value = 0.5334899687462222
$ws.Range("R4").Value = 4.801409718716
$ws.Range("S4").Value = 0.01687482416360171
$ws.Range("T4").Value = 0.01687482416360171
$ws.Range("I5").Value = 0.8420553458721338
$ws.Range("J5").Value = 0.8420553458721339
$ws.Range("M5").Value = 1.485259333333333
$ws.Range("N5").Value = 4.455778
$ws.Range("O5").Value = 0.3057455162066235
$ws.Range("P5").Value = 0.3057455162066235
$ws.Range("Q5").Value = 8.139312737658891
$ws.Range("R5").Value = 73.25381463893001
$ws.Range("S5").Value = 0.2574546463982225
$ws.Range("T5").Value = 0.2574546463982225
$ws.Range("I6").Value = 0.8420553458721338
$ws.Range("J6").Value = 0.8420553458721339
$ws.Range("O6").Value = 0.2805555239151429
$ws.Range("P6").Value = 0.2805555239151429
$ws.Range("S6").Value = 0.2362432787267034
$ws.Range("T6").Value = 0.2362432787267034
$ws.Range("I7").Value = 0.8420553458721338
$ws.Range("J7").Value = 0.8420553458721339
$ws.Range("O7").Value = 0.4136989598782336
$ws.Range("P7").Value = 0.4136989598782336
$ws.Range("S7").Value = 0.348357420747208
$ws.Range("T7").Value = 0.3483574207472081
$ws.Range("G8").Value = 0.7624369999999999
$ws.Range("I8").Value = 0.1171545487609863
$ws.Range("J8").Value = 0.1171545487609864
$ws.Range("M8").Value = 1.485259333333333
$ws.Range("N8").Value = 4.455778
$ws.Range("O8").Value = 0.3057455162066235
$ws.Range("P8").Value = 0.3057455162066235
$ws.Range("Q8").Value = 1.132416670328667
$ws.Range("R8").Value = 10.191750032958
$ws.Range("S8").Value = 0.03581947798688181
$ws.Range("T8").Value = 0.03581947798688181
$ws.Range("G9").Value = 0.7624369999999999
$ws.Range("I9").Value = 0.1171545487609863
$ws.Range("J9").Value = 0.1171545487609864
$ws.Range("O9").Value = 0.2805555239151429
$ws.Range("P9").Value = 0.2805555239151429
$ws.Range("Q9").Value = 1.039118271221333
$ws.Range("R9").Value = 9.352064440992001
$ws.Range("S9").Value = 0.03286835580668068
$ws.Range("T9").Value = 0.03286835580668069
$ws.Range("G10").Value = 0.7624369999999999
$ws.Range("I10").Value = 0.1171545487609863
$ws.Range("J10").Value = 0.1171545487609864
$ws.Range("O10").Value = 0.4136989598782336
$ws.Range("P10").Value = 0.4136989598782336
$ws.Range("S10").Value = 0.04846671496742385
$ws.Range("T10").Value = 0.04846671496742386
